# Atualizado por script em 27-11-2023 14:45
#
# 1) Rows 40 and 41 had their match details (columns F:V) swapped by
#    mistake; fix the order while keeping A:E (index/pais/torneio/
#    temporada/data_partida) exactly where they were.
# 2) A new match (row 53) is appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap F:V between row 40 and row 41
#    (NOTE: use .Value2 for the read side - .Value's getter isn't
#    reliable for multi-cell Range reads on this host)
# ---------------------------------------------------------------
$row40 = $ws.Range("F40:V40").Value2
$row41 = $ws.Range("F41:V41").Value2

$ws.Range("F40:V40").Value2 = $row41
$ws.Range("F41:V41").Value2 = $row40

# ---------------------------------------------------------------
# 2) Append new row 53, matching the look & feel (styles) of the
#    preceding data row (52): bordered/bold index column, date
#    formatted data_partida column, plain cells everywhere else.
# ---------------------------------------------------------------
$ws.Range("A52:V52").Copy()
$ws.Range("A53:V53").PasteSpecial(-4122)

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "moldova"
$ws.Range("C53").Value = "super-liga"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("E53").Value = 45257.5
$ws.Range("F53").Value = "Dacia Buiucani"
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = "Floresti"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 2.28
$ws.Range("K53").Value = "26/11/2023 01:12"
$ws.Range("L53").Value = 2.63
$ws.Range("M53").Value = "27/11/2023 11:59"
$ws.Range("N53").Value = 3.44
$ws.Range("O53").Value = "26/11/2023 01:12"
$ws.Range("P53").Value = 3.3
$ws.Range("Q53").Value = "27/11/2023 11:59"
$ws.Range("R53").Value = 2.62
$ws.Range("S53").Value = "26/11/2023 01:12"
$ws.Range("T53").Value = 2.37
$ws.Range("U53").Value = "27/11/2023 11:59"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/moldova/super-liga/dacia-buiucani-floresti/UHRAr4AR/"
